$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4116.346
$ws.Range("I62").Value = 3091.6667
$ws.Range("K62").Value = 3091.6667
$ws.Range("M62").Value = -2467.6667

$ws.Range("H65").Value = 4116.346
$ws.Range("I65").Value = 3091.6667
$ws.Range("K65").Value = 15458.3335
$ws.Range("M65").Value = -12338.3335

$ws.Range("H74").Value = 10752.889
$ws.Range("I74").Value = 10429.5
$ws.Range("K74").Value = 10429.5
$ws.Range("M74").Value = -9493.5

$ws.Range("H77").Value = 10752.889
$ws.Range("I77").Value = 10429.5
$ws.Range("K77").Value = 52147.5
$ws.Range("M77").Value = -47467.5

$ws.Range("H128").Value = 61250
$ws.Range("J128").Value = 61250
$ws.Range("L128").Value = 61250
$ws.Range("N128").Value = -71210

$ws.Range("H129").Value = 2515.95
$ws.Range("J129").Value = 2677.375
$ws.Range("L129").Value = 8032.125
$ws.Range("N129").Value = -18032.125

$ws.Range("H138").Value = 2486.9565
$ws.Range("J138").Value = 3255.0952
$ws.Range("L138").Value = 9765.285600000001
$ws.Range("N138").Value = -20045.2856

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 71433310
$ws.Range("I45").Value = 142858220
$ws.Range("J45").Value = 8394.714
$ws.Range("K45").Value = 142858220
$ws.Range("L45").Value = 8394.714
$ws.Range("M45").Value = -142857843
$ws.Range("N45").Value = -9148.714

$ws.Range("H63").Value = 2878.8
$ws.Range("J63").Value = 2897
$ws.Range("L63").Value = 2897
$ws.Range("N63").Value = -4269

$ws.Range("H66").Value = 2878.8
$ws.Range("J66").Value = 2897
$ws.Range("L66").Value = 14485
$ws.Range("N66").Value = -21349

$ws.Range("H76").Value = 9000
$ws.Range("J76").Value = 9000
$ws.Range("L76").Value = 9000
$ws.Range("N76").Value = -9676

$ws.Range("H79").Value = 9000
$ws.Range("J79").Value = 9000
$ws.Range("L79").Value = 9000
$ws.Range("N79").Value = -11340

$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5214.3335
$ws.Range("J86").Value = 9691.200000000001
$ws.Range("L86").Value = 9691.200000000001
$ws.Range("N86").Value = -11937.2

$ws.Range("H89").Value = 5214.3335
$ws.Range("J89").Value = 9691.200000000001
$ws.Range("L89").Value = 48456
$ws.Range("N89").Value = -59688

$ws.Range("H132").Value = 175050500
$ws.Range("J132").Value = 233378990
$ws.Range("L132").Value = 233378990
$ws.Range("N132").Value = -233389110

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 175.83333
$ws.Range("I7").Value = 66.22221999999999
$ws.Range("K7").Value = 66.22221999999999
$ws.Range("M7").Value = 46.77778000000001

$ws.Range("H132").Value = 8983.796
$ws.Range("I132").Value = 8762.24
$ws.Range("K132").Value = 26286.72
$ws.Range("M132").Value = -23756.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 160
$ws.Range("I41").Value = 150
$ws.Range("J41").Value = 166.66667
$ws.Range("K41").Value = 450
$ws.Range("L41").Value = 500.00001
$ws.Range("M41").Value = -112
$ws.Range("N41").Value = -1176.00001

$ws.Range("H59").Value = 1174.75
$ws.Range("I59").Value = 1174.75
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 3524.25
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -2984.25
$ws.Range("N59").ClearContents()

$ws.Range("H86").Value = 2052.4
$ws.Range("J86").Value = 3137
$ws.Range("L86").Value = 9411
$ws.Range("N86").Value = -11783

$ws.Range("H89").Value = 2052.4
$ws.Range("J89").Value = 3137
$ws.Range("L89").Value = 28233
$ws.Range("N89").Value = -40089

$ws.Range("H105").Value = 19360.666

$ws.Range("H120").Value = 4000.5
$ws.Range("I120").Value = 4000.5
$ws.Range("K120").Value = 12001.5
$ws.Range("M120").Value = -7163.5

$ws.Range("H125").Value = 18750
$ws.Range("J125").Value = 18750
$ws.Range("L125").Value = 56250
$ws.Range("N125").Value = -66090

$ws.Range("H133").Value = 1884.5
$ws.Range("I133").Value = 1884.5
$ws.Range("K133").Value = 5653.5
$ws.Range("M133").Value = -593.5

$ws.Range("H139").Value = 3098.9375
$ws.Range("I139").Value = 1646.4286
$ws.Range("J139").Value = 13266.5
$ws.Range("K139").Value = 4939.2858
$ws.Range("L139").Value = 39799.5
$ws.Range("M139").Value = 200.7142000000003
$ws.Range("N139").Value = -50079.5

$ws.Range("H140").Value = 3774.1667
$ws.Range("I140").Value = 2929
$ws.Range("J140").Value = 8000
$ws.Range("K140").Value = 8787
$ws.Range("L140").Value = 24000
$ws.Range("M140").Value = -3607
$ws.Range("N140").Value = -34360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3500
$ws.Range("I102").Value = 3666.6667
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 3666.6667
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = -2044.6667
$ws.Range("N102").Value = -6244

$ws.Range("H113").Value = 4353
$ws.Range("I113").Value = 3307.3333
$ws.Range("J113").Value = 7490
$ws.Range("K113").Value = 3307.3333
$ws.Range("L113").Value = 7490
$ws.Range("M113").Value = -1137.3333
$ws.Range("N113").Value = -11830

$ws.Range("H122").Value = 3363.4211
$ws.Range("I122").Value = 2641
$ws.Range("K122").Value = 7923
$ws.Range("M122").Value = -5473

$ws.Range("H136").Value = 41627.93
$ws.Range("J136").Value = 41627.93
$ws.Range("L136").Value = 124883.79
$ws.Range("N136").Value = -129983.79

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6375.125
$ws.Range("J22").Value = 9750.25
$ws.Range("L22").Value = 9750.25
$ws.Range("N22").Value = -10340.25

$ws.Range("H27").Value = 6375.125
$ws.Range("J27").Value = 9750.25
$ws.Range("L27").Value = 9750.25
$ws.Range("N27").Value = -9964.25

$ws.Range("H40").Value = 6894.1787
$ws.Range("I40").Value = 6834.619
$ws.Range("J40").Value = 7072.857
$ws.Range("K40").Value = 6834.619
$ws.Range("L40").Value = 7072.857
$ws.Range("M40").Value = -6698.619
$ws.Range("N40").Value = -7344.857

$ws.Range("H63").Value = 54900
$ws.Range("J63").Value = 54900
$ws.Range("L63").Value = 54900
$ws.Range("N63").Value = -56398

$ws.Range("H66").Value = 54900
$ws.Range("J66").Value = 54900
$ws.Range("L66").Value = 164700
$ws.Range("N66").Value = -172188

$ws.Range("H122").Value = 8640
$ws.Range("I122").Value = 7661.1113
$ws.Range("J122").Value = 10402
$ws.Range("K122").Value = 22983.3339
$ws.Range("L122").Value = 31206
$ws.Range("M122").Value = -20533.3339
$ws.Range("N122").Value = -36106

$ws.Range("H132").Value = 3054.8718
$ws.Range("J132").Value = 7310.3335
$ws.Range("L132").Value = 21931.0005
$ws.Range("N132").Value = -26991.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 17347.334
$ws.Range("J33").Value = 17347.334
$ws.Range("L33").Value = 17347.334
$ws.Range("N33").Value = -17847.334

$ws.Range("H36").Value = 17347.334
$ws.Range("J36").Value = 17347.334
$ws.Range("L36").Value = 17347.334
$ws.Range("N36").Value = -17847.334

$ws.Range("H70").Value = 28997
$ws.Range("I70").Value = 28997
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 28997
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -28682
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 28997
$ws.Range("I73").Value = 28997
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 28997
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -27905
$ws.Range("N73").ClearContents()

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H126").Value = 3978.3333
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3978.3333
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 11934.9999
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -16874.9999

$ws.Range("H136").Value = 3406.2942
$ws.Range("I136").Value = 2540.1333
$ws.Range("K136").Value = 7620.3999
$ws.Range("M136").Value = -5070.3999
